$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Job adverts by occupation" latest period value (C13)
# from "Nov 2024 (05/11/24)" to "Nov 2024 (07/02/25)"
$ws.Range("C13").Value = "Nov 2024 (07/02/25)"

# Update the selected cell to D13 (as reflected in the saved view state)
$ws.Range("D13").Select()
